$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns (AD, AE, AF)
# Reuse the existing header style (bold, centered, bordered) from AC1 by
# copying its formatting onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record values (Wins=71, Losses=91, Ties=0) for every data row
$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 91
    $ws.Cells.Item($r, 32).Value = 0
}
